$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.813.71'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '''3.338.37'
$ws.Range("E3").Value = '  +1.97%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''581.05'
$ws.Range("E5").Value = '  +0.38%  '
$ws.Range("D6").Value = '''176.13'
$ws.Range("E6").Value = '  +2.45%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +1.98%  '
$ws.Range("D9").Value = '''3.334.04'
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("E10").Value = '  +6.36%  '
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("D12").Value = '''46.95'
$ws.Range("E12").Value = '  +4.82%  '
$ws.Range("D13").Value = '''0.0000273'
$ws.Range("D14").Value = '''690.17'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '''3.877.80'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("E16").Value = '  +2.63%  '
$ws.Range("D17").Value = '''67.804.27'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").Value = '''3.327.11'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '''17.59'
$ws.Range("E20").Value = '  +2.56%  '
$ws.Range("D21").Value = '''11.07'
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("D23").Value = '''5.45'
$ws.Range("E23").Value = '  +4.50%  '
$ws.Range("D24").Value = '''16.98'
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = '''101.67'
$ws.Range("E25").Value = '  +3.44%  '
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("E27").Value = '  +2.46%  '
$ws.Range("D28").Value = '''9.52'
$ws.Range("E28").Value = '  +5.70%  '
$ws.Range("D29").Value = '''33.04'
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("E30").Value = '  +3.53%  '
$ws.Range("D31").Value = '''7.08'
$ws.Range("E31").Value = '  +7.74%  '
$ws.Range("D32").Value = '''569.07'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").Value = '''11.02'
$ws.Range("E33").Value = '  +2.29%  '
$ws.Range("E34").Value = '  +3.43%  '
$ws.Range("D35").Value = '''57.57'
$ws.Range("E35").Value = '  +4.11%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").Value = '''3.709.91'
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("D38").Value = '''3.30'
$ws.Range("E38").Value = '  +1.21%  '
$ws.Range("D39").Value = '''35.16'
$ws.Range("E39").Value = '  +12.60%  '
$ws.Range("D40").Value = '''0.135'
$ws.Range("E40").Value = '  +5.80%  '
$ws.Range("D41").Value = '''3.17'
$ws.Range("E41").Value = '  +7.34%  '
$ws.Range("E42").Value = '  +3.12%  '
$ws.Range("D43").Value = '''0.0₃0675'
$ws.Range("E43").Value = '  +3.48%  '
$ws.Range("E44").Value = '  +4.08%  '
$ws.Range("D45").Value = '''3.32'
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").Value = '''0.0409'
$ws.Range("E46").Value = '  +2.18%  '
$ws.Range("E47").Value = '  +5.85%  '
$ws.Range("E48").Value = '  +2.04%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").Value = '''132.17'
$ws.Range("E51").Value = '  +3.13%  '
